$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Nội dung câu trả lời"
$ws.Range("E6").Select()
